$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 369
$ws1.Range("F5").Value = 1789
$ws1.Range("F7").Value = 1466
$ws1.Range("F10").Value = 714
$ws1.Range("F11").Value = 13040
$ws1.Range("F12").Value = 12957
$ws1.Range("F18").Value = 615
$ws1.Range("F22").Value = 25
$ws1.Range("F24").Value = 154
$ws1.Range("F25").Value = 263
$ws1.Range("F26").Value = 719

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 28

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 369
$ws4.Range("F7").Value = 1789
$ws4.Range("F9").Value = 1466
$ws4.Range("F13").Value = 714
$ws4.Range("F14").Value = 13040
$ws4.Range("F15").Value = 12957
$ws4.Range("F21").Value = 615
$ws4.Range("F22").Value = 28
$ws4.Range("F27").Value = 25
$ws4.Range("F31").Value = 154
$ws4.Range("F32").Value = 263
$ws4.Range("F33").Value = 719
